$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-21 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-22 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("43-19=", $true, $false, $false, $false, $false, $true, 1, $false, "87-46=", 2) | Out-Null
$d.Content.Find.Execute("66-2=", $true, $false, $false, $false, $false, $true, 1, $false, "10+22=", 2) | Out-Null
$d.Content.Find.Execute("15+7=", $true, $false, $false, $false, $false, $true, 1, $false, "98-74=", 2) | Out-Null
$d.Content.Find.Execute("17-8=", $true, $false, $false, $false, $false, $true, 1, $false, "68+5=", 2) | Out-Null
$d.Content.Find.Execute("6+46=", $true, $false, $false, $false, $false, $true, 1, $false, "17+72=", 2) | Out-Null
$d.Content.Find.Execute("38+6=", $true, $false, $false, $false, $false, $true, 1, $false, "79-44=", 2) | Out-Null
$d.Content.Find.Execute("93-10=", $true, $false, $false, $false, $false, $true, 1, $false, "61+4=", 2) | Out-Null
$d.Content.Find.Execute("47+11=", $true, $false, $false, $false, $false, $true, 1, $false, "41+57=", 2) | Out-Null
$d.Content.Find.Execute("29-9=", $true, $false, $false, $false, $false, $true, 1, $false, "37+54=", 2) | Out-Null
$d.Content.Find.Execute("77-18=", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=", 2) | Out-Null
$d.Content.Find.Execute("63+30=", $true, $false, $false, $false, $false, $true, 1, $false, "42+21=", 2) | Out-Null
$d.Content.Find.Execute("64+1=", $true, $false, $false, $false, $false, $true, 1, $false, "43+18=", 2) | Out-Null
$d.Content.Find.Execute("67-33=", $true, $false, $false, $false, $false, $true, 1, $false, "48+4=", 2) | Out-Null
$d.Content.Find.Execute("39+8=", $true, $false, $false, $false, $false, $true, 1, $false, "66-35=", 2) | Out-Null
$d.Content.Find.Execute("56-5=", $true, $false, $false, $false, $false, $true, 1, $false, "2+47=", 2) | Out-Null
$d.Content.Find.Execute("77-12=", $true, $false, $false, $false, $false, $true, 1, $false, "90-76=", 2) | Out-Null
$d.Content.Find.Execute("78-60=", $true, $false, $false, $false, $false, $true, 1, $false, "79-6=", 2) | Out-Null
$d.Content.Find.Execute("27+44=", $true, $false, $false, $false, $false, $true, 1, $false, "42-2=", 2) | Out-Null
$d.Content.Find.Execute("58+30=", $true, $false, $false, $false, $false, $true, 1, $false, "87-28=", 2) | Out-Null
$d.Content.Find.Execute("3+50=", $true, $false, $false, $false, $false, $true, 1, $false, "42-41=", 2) | Out-Null
$d.Content.Find.Execute("73-41=", $true, $false, $false, $false, $false, $true, 1, $false, "94-83=", 2) | Out-Null
$d.Content.Find.Execute("2+14=", $true, $false, $false, $false, $false, $true, 1, $false, "41-12=", 2) | Out-Null
$d.Content.Find.Execute("53-32=", $true, $false, $false, $false, $false, $true, 1, $false, "0+10=", 2) | Out-Null
$d.Content.Find.Execute("12+70=", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=", 2) | Out-Null
$d.Content.Find.Execute("78+7=", $true, $false, $false, $false, $false, $true, 1, $false, "11-2=", 2) | Out-Null
$d.Content.Find.Execute("69-18=", $true, $false, $false, $false, $false, $true, 1, $false, "25+6=", 2) | Out-Null
$d.Content.Find.Execute("77-66=", $true, $false, $false, $false, $false, $true, 1, $false, "64+0=", 2) | Out-Null
$d.Content.Find.Execute("32-25=", $true, $false, $false, $false, $false, $true, 1, $false, "48+20=", 2) | Out-Null
$d.Content.Find.Execute("39+47=", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=", 2) | Out-Null
$d.Content.Find.Execute("88-72=", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=", 2) | Out-Null
$d.Content.Find.Execute("76-24=", $true, $false, $false, $false, $false, $true, 1, $false, "85-85=", 2) | Out-Null
$d.Content.Find.Execute("1+21=", $true, $false, $false, $false, $false, $true, 1, $false, "67-36=", 2) | Out-Null
$d.Content.Find.Execute("3+64=", $true, $false, $false, $false, $false, $true, 1, $false, "80+15=", 2) | Out-Null
$d.Content.Find.Execute("28+68=", $true, $false, $false, $false, $false, $true, 1, $false, "79-32=", 2) | Out-Null
$d.Content.Find.Execute("88-66=", $true, $false, $false, $false, $false, $true, 1, $false, "17+54=", 2) | Out-Null
$d.Content.Find.Execute("0+72=", $true, $false, $false, $false, $false, $true, 1, $false, "68+10=", 2) | Out-Null
$d.Content.Find.Execute("55+3=", $true, $false, $false, $false, $false, $true, 1, $false, "58+7=", 2) | Out-Null
$d.Content.Find.Execute("6+65=", $true, $false, $false, $false, $false, $true, 1, $false, "38-31=", 2) | Out-Null
$d.Content.Find.Execute("18+20=", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=", 2) | Out-Null
$d.Content.Find.Execute("57+26=", $true, $false, $false, $false, $false, $true, 1, $false, "99-7=", 2) | Out-Null
$d.Content.Find.Execute("87-1=", $true, $false, $false, $false, $false, $true, 1, $false, "98-34=", 2) | Out-Null
$d.Content.Find.Execute("41+15=", $true, $false, $false, $false, $false, $true, 1, $false, "19+15=", 2) | Out-Null
$d.Content.Find.Execute("43-2=", $true, $false, $false, $false, $false, $true, 1, $false, "16+67=", 2) | Out-Null
$d.Content.Find.Execute("55-6=", $true, $false, $false, $false, $false, $true, 1, $false, "81-9=", 2) | Out-Null
$d.Content.Find.Execute("70-66=", $true, $false, $false, $false, $false, $true, 1, $false, "43+42=", 2) | Out-Null
$d.Content.Find.Execute("70-20=", $true, $false, $false, $false, $false, $true, 1, $false, "83-43=", 2) | Out-Null
$d.Content.Find.Execute("35+18=", $true, $false, $false, $false, $false, $true, 1, $false, "67+0=", 2) | Out-Null
$d.Content.Find.Execute("33+49=", $true, $false, $false, $false, $false, $true, 1, $false, "13-9=", 2) | Out-Null
$d.Content.Find.Execute("20+52=", $true, $false, $false, $false, $false, $true, 1, $false, "43-3=", 2) | Out-Null
$d.Content.Find.Execute("63-27=", $true, $false, $false, $false, $false, $true, 1, $false, "7+15=", 2) | Out-Null
$d.Content.Find.Execute("26-25=", $true, $false, $false, $false, $false, $true, 1, $false, "9+69=", 2) | Out-Null
$d.Content.Find.Execute("23-5=", $true, $false, $false, $false, $false, $true, 1, $false, "82-34=", 2) | Out-Null
$d.Content.Find.Execute("6+23=", $true, $false, $false, $false, $false, $true, 1, $false, "5+90=", 2) | Out-Null
$d.Content.Find.Execute("78+19=", $true, $false, $false, $false, $false, $true, 1, $false, "77-34=", 2) | Out-Null
$d.Content.Find.Execute("4+82=", $true, $false, $false, $false, $false, $true, 1, $false, "17+1=", 2) | Out-Null
$d.Content.Find.Execute("14+80=", $true, $false, $false, $false, $false, $true, 1, $false, "73-16=", 2) | Out-Null
$d.Content.Find.Execute("77+6=", $true, $false, $false, $false, $false, $true, 1, $false, "78-8=", 2) | Out-Null
$d.Content.Find.Execute("26+8=", $true, $false, $false, $false, $false, $true, 1, $false, "25+3=", 2) | Out-Null
$d.Content.Find.Execute("36+0=", $true, $false, $false, $false, $false, $true, 1, $false, "74-68=", 2) | Out-Null
$d.Content.Find.Execute("51-26=", $true, $false, $false, $false, $false, $true, 1, $false, "67+0=", 2) | Out-Null
$d.Content.Find.Execute("33+41=", $true, $false, $false, $false, $false, $true, 1, $false, "1+94=", 2) | Out-Null
$d.Content.Find.Execute("74+19=", $true, $false, $false, $false, $false, $true, 1, $false, "36-30=", 2) | Out-Null
$d.Content.Find.Execute("3+25=", $true, $false, $false, $false, $false, $true, 1, $false, "93-84=", 2) | Out-Null
$d.Content.Find.Execute("29+62=", $true, $false, $false, $false, $false, $true, 1, $false, "42-5=", 2) | Out-Null
$d.Content.Find.Execute("54-5=", $true, $false, $false, $false, $false, $true, 1, $false, "57+3=", 2) | Out-Null
$d.Content.Find.Execute("0+70=", $true, $false, $false, $false, $false, $true, 1, $false, "31+1=", 2) | Out-Null
$d.Content.Find.Execute("61+28=", $true, $false, $false, $false, $false, $true, 1, $false, "2+59=", 2) | Out-Null
$d.Content.Find.Execute("14+79=", $true, $false, $false, $false, $false, $true, 1, $false, "43+5=", 2) | Out-Null
$d.Content.Find.Execute("86-75=", $true, $false, $false, $false, $false, $true, 1, $false, "61-45=", 2) | Out-Null
$d.Content.Find.Execute("91-72=", $true, $false, $false, $false, $false, $true, 1, $false, "69-58=", 2) | Out-Null
$d.Content.Find.Execute("38+58=", $true, $false, $false, $false, $false, $true, 1, $false, "69+29=", 2) | Out-Null
$d.Content.Find.Execute("58+36=", $true, $false, $false, $false, $false, $true, 1, $false, "10+89=", 2) | Out-Null
$d.Content.Find.Execute("65-51=", $true, $false, $false, $false, $false, $true, 1, $false, "46-5=", 2) | Out-Null
$d.Content.Find.Execute("10+37=", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$d.Content.Find.Execute("30+63=", $true, $false, $false, $false, $false, $true, 1, $false, "69+0=", 2) | Out-Null
$d.Content.Find.Execute("86-23=", $true, $false, $false, $false, $false, $true, 1, $false, "11-10=", 2) | Out-Null
$d.Content.Find.Execute("99-80=", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=", 2) | Out-Null
$d.Content.Find.Execute("65+0=", $true, $false, $false, $false, $false, $true, 1, $false, "17+54=", 2) | Out-Null
$d.Content.Find.Execute("84-36=", $true, $false, $false, $false, $false, $true, 1, $false, "73+23=", 2) | Out-Null
$d.Content.Find.Execute("41+44=", $true, $false, $false, $false, $false, $true, 1, $false, "88-84=", 2) | Out-Null
$d.Content.Find.Execute("56+3=", $true, $false, $false, $false, $false, $true, 1, $false, "75-24=", 2) | Out-Null
$d.Content.Find.Execute("24+62=", $true, $false, $false, $false, $false, $true, 1, $false, "23+63=", 2) | Out-Null
$d.Content.Find.Execute("57+23=", $true, $false, $false, $false, $false, $true, 1, $false, "16+14=", 2) | Out-Null
$d.Content.Find.Execute("81-53=", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=", 2) | Out-Null
$d.Content.Find.Execute("93-87=", $true, $false, $false, $false, $false, $true, 1, $false, "42+47=", 2) | Out-Null
$d.Content.Find.Execute("36+17=", $true, $false, $false, $false, $false, $true, 1, $false, "60-15=", 2) | Out-Null
$d.Content.Find.Execute("22-8=", $true, $false, $false, $false, $false, $true, 1, $false, "95-85=", 2) | Out-Null
$d.Content.Find.Execute("60-13=", $true, $false, $false, $false, $false, $true, 1, $false, "27-8=", 2) | Out-Null
$d.Content.Find.Execute("26+0=", $true, $false, $false, $false, $false, $true, 1, $false, "97-9=", 2) | Out-Null
$d.Content.Find.Execute("61+34=", $true, $false, $false, $false, $false, $true, 1, $false, "31+65=", 2) | Out-Null
$d.Content.Find.Execute("91-62=", $true, $false, $false, $false, $false, $true, 1, $false, "15+68=", 2) | Out-Null
$d.Content.Find.Execute("4-4=", $true, $false, $false, $false, $false, $true, 1, $false, "83-40=", 2) | Out-Null
$d.Content.Find.Execute("62-8=", $true, $false, $false, $false, $false, $true, 1, $false, "62+4=", 2) | Out-Null
$d.Content.Find.Execute("21+6=", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=", 2) | Out-Null
$d.Content.Find.Execute("58+22=", $true, $false, $false, $false, $false, $true, 1, $false, "2+75=", 2) | Out-Null
$d.Content.Find.Execute("41-2=", $true, $false, $false, $false, $false, $true, 1, $false, "59-6=", 2) | Out-Null
$d.Content.Find.Execute("1+41=", $true, $false, $false, $false, $false, $true, 1, $false, "29+55=", 2) | Out-Null
$d.Content.Find.Execute("54-29=", $true, $false, $false, $false, $false, $true, 1, $false, "81-0=", 2) | Out-Null
$d.Content.Find.Execute("89-38=", $true, $false, $false, $false, $false, $true, 1, $false, "92-92=", 2) | Out-Null
$d.Content.Find.Execute("13+83=", $true, $false, $false, $false, $false, $true, 1, $false, "51-33=", 2) | Out-Null
